# Weekly fruit/vegetable price update: insert a new "Cilantro" record dated
# 2021-12-21 (serial 44551) right after the existing 2021-07-14 row (row 54),
# pushing the subsequent rows down by one. The new record repeats the most
# recent week's (2021-12-06) price figures, per the "Fruta / hortaliza,
# semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 54; existing rows 54-58 shift down to 55-59.
$ws.Rows("54:54").Insert()

$ws.Cells.Item(54, 1).Value  = 1
$ws.Cells.Item(54, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(54, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(54, 4).Value  = 44551
$ws.Cells.Item(54, 5).Value  = 15
$ws.Cells.Item(54, 6).Value  = 100112040
$ws.Cells.Item(54, 7).Value  = 'Cilantro'
$ws.Cells.Item(54, 8).Value  = 'Sin especificar'
$ws.Cells.Item(54, 9).Value  = 'Primera'
$ws.Cells.Item(54, 10).Value = 300
$ws.Cells.Item(54, 11).Value = 2800
$ws.Cells.Item(54, 12).Value = 3000
$ws.Cells.Item(54, 13).Value = 2900
$ws.Cells.Item(54, 14).Value = '$/atado 1,5 a 2 kilos'
$ws.Cells.Item(54, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(54, 16).Value = 1450
$ws.Cells.Item(54, 17).Value = 2
$ws.Cells.Item(54, 18).Value = 'Hortaliza'
